$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update B2 value (186 -> 342)
$ws.Range("B2").Value = 342

# Update A3 value (2 -> 1) and B3 value (156 -> 127)
$ws.Range("A3").Value = 1
$ws.Range("B3").Value = 127

# Remove old row 4 contents (A4=1, B4=127) entirely so the used range shrinks to A1:B3
$ws.Range("A4:B4").Delete()
